$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fill in the previously missing surname on row 3 (D3) and
# the previously missing name on row 4 (E4).
$ws.Range("D3").Value = "Фамилия 2"
$ws.Range("E4").Value = "Имя 3"

# Move the active selection to H10, matching the saved view state.
$ws.Range("H10").Select()
